$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the "Number" column to be stored as text so purely-numeric phone
# numbers are not reinterpreted as numbers (they must stay shared strings,
# matching the rest of column B).
$ws.Range("B355:B358").NumberFormat = "@"

# Row 355
$ws.Range("A355").Value = "2026-02-18 16:43:43"
$ws.Range("B355").Value = "237674895525"
$ws.Range("C355").Value = "NFOR EPSE FOMUNGUM ASSUMPTA MUMBEB"
$ws.Range("D355").Value = 9524

# Row 356
$ws.Range("A356").Value = "2026-02-18 14:01:57"
$ws.Range("B356").Value = "237683454307"
$ws.Range("C356").Value = "NGO MBOMNDA JULIENNE ELVA CHIC MOBILE SARL"
$ws.Range("D356").Value = 4717

# Row 357
$ws.Range("A357").Value = "2026-02-18 14:21:23"
$ws.Range("B357").Value = "237671104974"
$ws.Range("C357").Value = "MFS AM MARCHÉ BEEDI"
$ws.Range("D357").Value = 0

# Row 358
$ws.Range("A358").Value = "2026-02-18 12:35:56"
$ws.Range("B358").Value = "237671694408"
$ws.Range("C358").Value = "VARELLE NKENGUE BILLE"
$ws.Range("D358").Value = 128405

# Restore the default (General) look for column B so no extra visible
# formatting sticks to the new cells, while keeping the text storage type.
$ws.Range("B355:B358").ClearFormats()
